$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.963.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.63%  "
$ws.Range("D3").Value = "'4.020.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.38%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'535.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.80%  "
$ws.Range("D6").Value = "'153.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.88%  "
$ws.Range("D7").Value = "'0.693"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +14.22%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.752"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.34%  "
$ws.Range("D10").Value = "'0.173"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").Value = "'0.0000326"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").Value = "'47.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +15.53%  "
$ws.Range("D13").Value = "'10.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.58%  "
$ws.Range("D14").Value = "'4.651.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.26%  "
$ws.Range("D15").Value = "'4.016.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.00%  "
$ws.Range("D16").Value = "'14.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").Value = "'20.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "'71.722.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.49%  "
$ws.Range("D21").Value = "'432.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.18%  "
$ws.Range("D22").Value = "'98.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.64%  "
$ws.Range("D23").Value = "'3.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.39%  "
$ws.Range("D24").Value = "'14.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.39%  "
$ws.Range("D25").Value = "'4.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.08%  "
$ws.Range("D26").Value = "'11.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.74%  "
$ws.Range("D27").Value = "'10.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.79%  "
$ws.Range("D28").Value = "'3.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +31.27%  "
$ws.Range("D29").Value = "'5.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.89%  "
$ws.Range("D30").Value = "'36.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.27%  "
$ws.Range("D31").Value = "'13.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("E32").Value = "  +5.31%  "
$ws.Range("D33").Value = "'681.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("D34").Value = "'6.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("D35").Value = "'65.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").Value = "'42.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.67%  "
$ws.Range("D37").Value = "'0.426"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.14%  "
$ws.Range("D38").Value = "'3.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.19%  "
$ws.Range("D39").Value = "'0.155"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.17%  "
$ws.Range("D40").Value = "'0.0₃0829"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("D41").Value = "'3.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "'0.996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").Value = "'0.0487"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.92%  "
$ws.Range("D45").Value = "'0.151"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.58%  "
$ws.Range("E46").Value = "  -7.70%  "
$ws.Range("D47").Value = "'9.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.75%  "
$ws.Range("E48").Value = "  -6.74%  "
$ws.Range("D49").Value = "'3.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").Value = "'3.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").Value = "'144.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.92%  "
